# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last refreshed" timestamp on row 1
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 15:05"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1837830
$ws.Range("C4").Value = 660
$ws.Range("E4").Value = 1131735
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = 106220

# Row 10: Peru
$ws.Range("B10").Value = 191605
$ws.Range("C10").Value = 996
$ws.Range("D10").Value = 92045
$ws.Range("E10").Value = 94145
$ws.Range("G10").Value = 7
$ws.Range("H10").Value = 5415

# Row 12: Mexico
$ws.Range("B12").Value = 183515
$ws.Range("C12").Value = 21
$ws.Range("E12").Value = 9010

# Row 19: Paises Bajos
$ws.Range("B19").Value = 87142
$ws.Range("C19").Value = 1881
$ws.Range("D19").Value = 64306
$ws.Range("E19").Value = 22311
$ws.Range("G19").Value = 22
$ws.Range("H19").Value = 525

# Rows 31/32 swap order (Portugal overtakes Sudafrica): update labels and values
$ws.Range("A31").Value = "Portugal"
$ws.Range("B31").Value = 32700
$ws.Range("C31").Value = 200
$ws.Range("D31").Value = 19552
$ws.Range("E31").Value = 11724
$ws.Range("G31").Value = 14
$ws.Range("H31").Value = 1424

$ws.Range("A32").Value = "Sudafrica"
$ws.Range("B32").Value = 32683
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 16809
$ws.Range("E32").Value = 15191
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 683

# Row 79: Uzbekistan
$ws.Range("B79").Value = 3673
$ws.Range("C79").Value = 50
$ws.Range("E79").Value = 799

# Row 102: Sri Lanka
$ws.Range("B102").Value = 1635
$ws.Range("C102").Value = 2
$ws.Range("E102").Value = 813

# Row 111: Libano
$ws.Range("B111").Value = 1233
$ws.Range("C111").Value = 13
$ws.Range("D111").Value = 715
$ws.Range("E111").Value = 491

# Row 123: Sierra Leona
$ws.Range("B123").Value = 865
$ws.Range("C123").Value = 4
$ws.Range("D123").Value = 475
$ws.Range("E123").Value = 344

# Row 171: Siria
$ws.Range("B171").Value = 123
$ws.Range("C171").Value = 1
$ws.Range("E171").Value = 72
